$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.325.51"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.839.90"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'239.12"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'0.6248"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'0.07372"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").Value = "'0.2883"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'24.70"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").Value = "'0.07724"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "1.842.89"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "'4.951"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "'0.6633"
$ws.Range("E15").Value = "  -2.81%  "
$ws.Range("D16").Value = "'81.24"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "'6.250"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "29.294.83"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "'233.32"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'157.38"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Value = "'8.425"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("D26").Value = "'0.1336"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "'0.07134"
$ws.Range("E28").Value = "  +7.96%  "
$ws.Range("D29").Value = "'1.487"
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").Value = "'4.031"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").Value = "'1.150"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").Value = "'1.809"
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("D35").Value = "'0.6947"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("E36").Value = "  +0.54%  "
$ws.Range("D37").Value = "'0.01826"
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("D38").Value = "'2.781"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").Value = "1.231.59"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("D40").Value = "'6.780"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'0.9437"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'101.11"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "1.985.14"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("D45").Value = "'65.15"
$ws.Range("E45").Value = "  -1.77%  "
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").Value = "'6.936"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").Value = "'1.680"
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("D49").Value = "'8.909"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").Value = "'0.1130"
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("D51").Value = "'0.3868"
$ws.Range("E51").Value = "  -2.23%  "
